$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worlds-like test")

# -----------------------------------------------------------------
# 1) Move the "AVG X/Y ERR" block from rows 22/23/25-28 up to
#    rows 17/18/20-23, and update K12/L12 to reference the new
#    location of the average (row 18 instead of row 23).
# -----------------------------------------------------------------
$ws.Range("K22:L23").ClearContents()
$ws.Range("K25:L28").ClearContents()

$ws.Range("K17").Value = "AVG X ERR"
$ws.Range("L17").Value = "AVG Y ERR"

$ws.Range("K18").Formula = "=AVERAGE(K2:K5, K8)"
$ws.Range("L18").Formula = "=AVERAGE(L2:L5, L8)"

$ws.Range("K20").Formula = "=K2-`$K`$18"
$ws.Range("L20").Formula = "=L2-`$L`$18"
$ws.Range("K21").Formula = "=K3-`$K`$18"
$ws.Range("L21").Formula = "=L3-`$L`$18"
$ws.Range("K22").Formula = "=K4-`$K`$18"
$ws.Range("L22").Formula = "=L4-`$L`$18"
$ws.Range("K23").Formula = "=K5-`$K`$18"
$ws.Range("L23").Formula = "=L5-`$L`$18"

$ws.Range("K12").Formula = "=K8-K18"
$ws.Range("L12").Formula = "=L8-L18"

# -----------------------------------------------------------------
# 2) Add the new "R1..R4" test-coordinate rows (29-34).
# -----------------------------------------------------------------
$ws.Range("H29").Formula = "=3/8"

# Copy the number formatting used by the other coordinate rows
# down onto the new rows first, then fill in the actual values.
$ws.Range("B2:D2").Copy()
$ws.Range("B30:D33").PasteSpecial(-4122)

$ws.Range("E8:F8").Copy()
$ws.Range("E30:F33").PasteSpecial(-4122)
$ws.Range("E30:F33").Font.ThemeColor = 1

$ws.Range("A6:M6").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A30").Value = "R1"
$ws.Range("B30").Value = 48
$ws.Range("C30").Value = -44
$ws.Range("E30").Value = 23.875
$ws.Range("F30").Value = 47.625
$ws.Range("K30").Formula = "=E30 - ABS(B30)"
$ws.Range("L30").Formula = "=F30 - ABS(C30)"
$ws.Range("M30").Formula = "=G30 - ABS(D30)"

$ws.Range("A31").Value = "R2"
$ws.Range("B31").Value = 23.7
$ws.Range("C31").Value = -44.5
$ws.Range("E31").Value = 47.625
$ws.Range("F31").Value = 47.625
$ws.Range("K31").Formula = "=E31 - ABS(B31)"
$ws.Range("L31").Formula = "=F31 - ABS(C31)"
$ws.Range("M31").Formula = "=G31 - ABS(D31)"

$ws.Range("A32").Value = "R3"
$ws.Range("B32").Value = 47.4
$ws.Range("C32").Value = -19.9
$ws.Range("E32").Value = 23.875
$ws.Range("F32").Value = 23.375
$ws.Range("K32").Formula = "=E32 - ABS(B32)"
$ws.Range("L32").Formula = "=F32 - ABS(C32)"
$ws.Range("M32").Formula = "=G32 - ABS(D32)"

$ws.Range("A33").Value = "R4"
$ws.Range("B33").Value = 23.6
$ws.Range("C33").Value = -20.3
$ws.Range("E33").Value = 47.625
$ws.Range("F33").Value = 23.375
$ws.Range("K33").Formula = "=E33 - ABS(B33)"
$ws.Range("L33").Formula = "=F33 - ABS(C33)"
$ws.Range("M33").Formula = "=G33 - ABS(D33)"

# -----------------------------------------------------------------
# 3) Update the sheet's view state (zoom + selected cell).
# -----------------------------------------------------------------
[void]$ws.Range("E30").Select()
$ws.Application.ActiveWindow.Zoom = 83

Write-Output "done"
